# Auto-generated Excel COM-interop edit script
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("LP1912")
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws3 = $wb.Worksheets.Item("6203-6173")

# --- LP1912: header text updates ---
$ws1.Cells.Item(2, 1).Value = "Última actualización: 08:47:26"
$ws1.Cells.Item(3, 1).Value = "Total filas: 128"

# --- LP1912-215: header text updates ---
$ws2.Cells.Item(2, 1).Value = "Última actualización: 08:47:26"
$ws2.Cells.Item(3, 1).Value = "Total filas: 20"

# --- 6203-6173: header text updates ---
$ws3.Cells.Item(2, 1).Value = "Última actualización: 08:47:26"
$ws3.Cells.Item(3, 1).Value = "Total filas: 21"

# --- LP1912: data rows (row 6 onward) ---
$ws1Data = @(
    @("04:45:48", "04:45", "215A_EL PATO", 0, "LP1912"),
    @("04:21:09", "04:46", "215A_EL PATO", 25, "LP1912"),
    @("04:45:48", "04:53", "11_ETCHEVERRY", 8, "LP1912"),
    @("04:56:30", "05:16", "17_ROMERO", 20, "LP1912"),
    @("04:56:30", "05:22", "23_HERNANDEZ", 26, "LP1912"),
    @("05:24:16", "05:25", "23_HERNANDEZ", 1, "LP1912"),
    @("04:56:30", "05:34", "215B_EL PATO", 38, "LP1912"),
    @("05:24:16", "05:35", "215B_EL PATO", 11, "LP1912"),
    @("05:24:16", "05:46", "15_ABASTO", 22, "LP1912"),
    @("05:24:16", "05:54", "10_OLMOS", 30, "LP1912"),
    @("05:55:02", "05:55", "10_OLMOS", 0, "LP1912"),
    @("05:24:16", "06:04", "16_SANTA ANA", 40, "LP1912"),
    @("05:24:16", "06:11", "215A_EL PATO", 47, "LP1912"),
    @("05:55:02", "06:12", "215A_EL PATO", 17, "LP1912"),
    @("05:55:02", "06:14", "225_HARAS DEL SUR", 19, "LP1912"),
    @("04:56:30", "06:18", "16_SANTA ANA", 82, "LP1912"),
    @("05:55:02", "06:21", "26_HERNANDEZ", 26, "LP1912"),
    @("04:45:48", "06:24", "16_SANTA ANA", 99, "LP1912"),
    @("05:55:02", "06:27", "23_HERNANDEZ", 32, "LP1912"),
    @("06:25:28", "06:29", "86_EST CHICA-ESC AGRARIA", 4, "LP1912"),
    @("06:25:28", "06:30", "23_HERNANDEZ", 5, "LP1912"),
    @("05:55:02", "06:30", "86_EST CHICA-ESC AGRARIA", 35, "LP1912"),
    @("05:55:02", "06:31", "16_SANTA ANA", 36, "LP1912"),
    @("06:25:28", "06:44", "225_C ROCA-H SUR", 19, "LP1912"),
    @("06:25:28", "06:46", "215C_EL PATO", 21, "LP1912"),
    @("05:55:02", "06:47", "215C_EL PATO", 52, "LP1912"),
    @("06:54:06", "06:55", "215C_EL PATO", 1, "LP1912"),
    @("06:54:06", "06:55", "14_ABASTO", 1, "LP1912"),
    @("06:25:28", "06:59", "14_ABASTO", 34, "LP1912"),
    @("05:55:02", "07:00", "14_ABASTO", 65, "LP1912"),
    @("06:54:06", "07:01", "16_SANTA ANA", 7, "LP1912"),
    @("06:25:28", "07:05", "15_ABASTO", 40, "LP1912"),
    @("06:54:06", "07:05", "23_HERNANDEZ", 11, "LP1912"),
    @("06:54:06", "07:07", "225_GOMEZ", 13, "LP1912"),
    @("06:54:06", "07:07", "15_ABASTO", 13, "LP1912"),
    @("06:25:28", "07:11", "215A_EL PATO", 46, "LP1912"),
    @("06:54:06", "07:12", "215A_EL PATO", 18, "LP1912"),
    @("06:25:28", "07:15", "11_ETCHEVERRY", 50, "LP1912"),
    @("06:54:06", "07:16", "11_ETCHEVERRY", 22, "LP1912"),
    @("06:54:06", "07:17", "16_SANTA ANA", 23, "LP1912"),
    @("07:17:59", "07:20", "26_HERNANDEZ", 3, "LP1912"),
    @("06:54:06", "07:21", "26_HERNANDEZ", 27, "LP1912"),
    @("06:54:06", "07:23", "10_OLMOS", 29, "LP1912"),
    @("07:17:59", "07:31", "84_COLONIA URQUIZA-ESC 49", 14, "LP1912"),
    @("07:17:59", "07:31", "16_SANTA ANA", 14, "LP1912"),
    @("07:17:59", "07:31", "11_ETCHEVERRY", 14, "LP1912"),
    @("06:54:06", "07:32", "11_ETCHEVERRY", 38, "LP1912"),
    @("06:54:06", "07:32", "84_COLONIA URQUIZA-ESC 49", 38, "LP1912"),
    @("05:55:02", "07:32", "16_SANTA ANA", 97, "LP1912"),
    @("07:17:59", "07:35", "23_HERNANDEZ", 18, "LP1912"),
    @("07:17:59", "07:36", "27_EL RETIRO", 19, "LP1912"),
    @("06:54:06", "07:37", "27_EL RETIRO", 43, "LP1912"),
    @("07:17:59", "07:38", "10_OLMOS", 21, "LP1912"),
    @("06:54:06", "07:39", "10_OLMOS", 45, "LP1912"),
    @("07:17:59", "07:46", "16_SANTA ANA", 29, "LP1912"),
    @("07:17:59", "07:47", "14_ABASTO", 30, "LP1912"),
    @("07:48:05", "07:48", "14_ABASTO", 0, "LP1912"),
    @("07:17:59", "07:51", "215D_EL PATO", 34, "LP1912"),
    @("07:48:05", "07:52", "215D_EL PATO", 4, "LP1912"),
    @("07:17:59", "07:59", "23_HERNANDEZ", 42, "LP1912"),
    @("07:48:05", "08:02", "23_HERNANDEZ", 14, "LP1912"),
    @("08:01:10", "08:03", "23_HERNANDEZ", 2, "LP1912"),
    @("07:17:59", "08:03", "11_ETCHEVERRY", 46, "LP1912"),
    @("08:01:10", "08:04", "11_ETCHEVERRY", 3, "LP1912"),
    @("06:54:06", "08:06", "23_HERNANDEZ", 72, "LP1912"),
    @("08:01:10", "08:11", "16_SANTA ANA", 10, "LP1912"),
    @("07:17:59", "08:11", "15_ABASTO", 54, "LP1912"),
    @("08:01:10", "08:12", "15_ABASTO", 11, "LP1912"),
    @("07:17:59", "08:12", "10_OLMOS", 55, "LP1912"),
    @("08:01:10", "08:13", "10_OLMOS", 12, "LP1912"),
    @("07:17:59", "08:20", "26_HERNANDEZ", 63, "LP1912"),
    @("06:54:06", "08:21", "26_HERNANDEZ", 87, "LP1912"),
    @("07:17:59", "08:22", "215B_EL PATO", 65, "LP1912"),
    @("07:17:59", "08:22", "16_P MOR-SANTA ANA", 65, "LP1912"),
    @("08:01:10", "08:23", "215B_EL PATO", 22, "LP1912"),
    @("08:01:10", "08:23", "16_P MOR-SANTA ANA", 22, "LP1912"),
    @("08:01:10", "08:24", "14_ABASTO", 23, "LP1912"),
    @("07:17:59", "08:26", "84_COLONIA URQUIZA-ESC 49", 69, "LP1912"),
    @("08:01:10", "08:27", "84_COLONIA URQUIZA-ESC 49", 26, "LP1912"),
    @("07:48:05", "08:30", "23_HERNANDEZ", 42, "LP1912"),
    @("08:01:10", "08:33", "10_OLMOS", 32, "LP1912"),
    @("08:31:01", "08:35", "23_HERNANDEZ", 4, "LP1912"),
    @("07:48:05", "08:37", "26_HERNANDEZ", 49, "LP1912"),
    @("08:01:10", "08:40", "16_SANTA ANA", 39, "LP1912"),
    @("07:17:59", "08:41", "81_EL PELIGRO", 84, "LP1912"),
    @("08:31:01", "08:42", "81_EL PELIGRO", 11, "LP1912"),
    @("07:17:59", "08:43", "14_ABASTO", 86, "LP1912"),
    @("08:31:01", "08:44", "16_SANTA ANA", 13, "LP1912"),
    @("08:31:01", "08:44", "14_ABASTO", 13, "LP1912"),
    @("08:47:26", "08:47", "16_SANTA ANA", 0, "LP1912"),
    @("08:01:10", "08:49", "26_HERNANDEZ", 48, "LP1912"),
    @("08:47:26", "08:53", "10_OLMOS", 6, "LP1912"),
    @("07:17:59", "08:53", "17_ROMERO", 96, "LP1912"),
    @("08:47:26", "08:54", "17_ROMERO", 7, "LP1912"),
    @("08:47:26", "09:01", "215A_EL PATO", 14, "LP1912"),
    @("08:31:01", "09:02", "215A_EL PATO", 31, "LP1912"),
    @("08:47:26", "09:03", "11_ETCHEVERRY", 16, "LP1912"),
    @("08:31:01", "09:04", "11_ETCHEVERRY", 33, "LP1912"),
    @("08:47:26", "09:05", "23_HERNANDEZ", 18, "LP1912"),
    @("08:47:26", "09:10", "16_P MOR-SANTA ANA", 23, "LP1912"),
    @("08:31:01", "09:11", "16_P MOR-SANTA ANA", 40, "LP1912"),
    @("08:47:26", "09:16", "27_EL RETIRO", 29, "LP1912"),
    @("08:31:01", "09:17", "27_EL RETIRO", 46, "LP1912"),
    @("08:31:01", "09:19", "26_HERNANDEZ", 48, "LP1912"),
    @("08:47:26", "09:21", "26_HERNANDEZ", 34, "LP1912"),
    @("08:47:26", "09:22", "16_SANTA ANA", 35, "LP1912"),
    @("08:31:01", "09:23", "17_ROMERO", 52, "LP1912"),
    @("08:31:01", "09:23", "16_SANTA ANA", 52, "LP1912"),
    @("08:47:26", "09:23", "11_ETCHEVERRY", 36, "LP1912"),
    @("08:31:01", "09:24", "11_ETCHEVERRY", 53, "LP1912"),
    @("07:48:05", "09:32", "15_ABASTO", 104, "LP1912"),
    @("08:47:26", "09:33", "10_OLMOS", 46, "LP1912"),
    @("07:48:05", "09:34", "16_SANTA ANA", 106, "LP1912"),
    @("08:47:26", "09:34", "23_HERNANDEZ", 47, "LP1912"),
    @("08:47:26", "09:35", "26_HERNANDEZ", 48, "LP1912"),
    @("08:47:26", "09:42", "215C_EL PATO", 55, "LP1912"),
    @("08:47:26", "09:43", "14_ABASTO", 56, "LP1912"),
    @("08:31:01", "09:44", "14_ABASTO", 73, "LP1912"),
    @("08:47:26", "09:52", "15_ABASTO", 65, "LP1912"),
    @("08:47:26", "09:53", "10_OLMOS", 66, "LP1912"),
    @("08:31:01", "09:56", "10_OLMOS", 85, "LP1912"),
    @("08:47:26", "10:10", "16_P MOR-SANTA ANA", 83, "LP1912"),
    @("08:31:01", "10:11", "16_P MOR-SANTA ANA", 100, "LP1912"),
    @("08:47:26", "10:12", "15_ABASTO", 85, "LP1912"),
    @("08:47:26", "10:26", "215A_EL PATO", 99, "LP1912"),
    @("08:31:01", "10:27", "215A_EL PATO", 116, "LP1912"),
    @("08:47:26", "10:42", "17_ROMERO", 115, "LP1912"),
    @("08:47:26", "10:43", "14_ABASTO", 116, "LP1912")
)
for ($i = 0; $i -lt $ws1Data.Length; $i++) {
    $r = 6 + $i
    for ($c = 0; $c -lt 5; $c++) {
        $ws1.Cells.Item($r, $c + 1).Value = $ws1Data[$i][$c]
    }
}

# --- LP1912-215: data rows (row 6 onward) ---
$ws2Data = @(
    @("04:45:48", "04:45", "215A_EL PATO", 0, "LP1912"),
    @("04:21:09", "04:46", "215A_EL PATO", 25, "LP1912"),
    @("04:56:30", "05:34", "215B_EL PATO", 38, "LP1912"),
    @("05:24:16", "05:35", "215B_EL PATO", 11, "LP1912"),
    @("05:24:16", "06:11", "215A_EL PATO", 47, "LP1912"),
    @("05:55:02", "06:12", "215A_EL PATO", 17, "LP1912"),
    @("06:25:28", "06:46", "215C_EL PATO", 21, "LP1912"),
    @("05:55:02", "06:47", "215C_EL PATO", 52, "LP1912"),
    @("06:54:06", "06:55", "215C_EL PATO", 1, "LP1912"),
    @("06:25:28", "07:11", "215A_EL PATO", 46, "LP1912"),
    @("06:54:06", "07:12", "215A_EL PATO", 18, "LP1912"),
    @("07:17:59", "07:51", "215D_EL PATO", 34, "LP1912"),
    @("07:48:05", "07:52", "215D_EL PATO", 4, "LP1912"),
    @("07:17:59", "08:22", "215B_EL PATO", 65, "LP1912"),
    @("08:01:10", "08:23", "215B_EL PATO", 22, "LP1912"),
    @("08:47:26", "09:01", "215A_EL PATO", 14, "LP1912"),
    @("08:31:01", "09:02", "215A_EL PATO", 31, "LP1912"),
    @("08:47:26", "09:42", "215C_EL PATO", 55, "LP1912"),
    @("08:47:26", "10:26", "215A_EL PATO", 99, "LP1912"),
    @("08:31:01", "10:27", "215A_EL PATO", 116, "LP1912")
)
for ($i = 0; $i -lt $ws2Data.Length; $i++) {
    $r = 6 + $i
    for ($c = 0; $c -lt 5; $c++) {
        $ws2.Cells.Item($r, $c + 1).Value = $ws2Data[$i][$c]
    }
}

# --- 6203-6173: data rows (row 6 onward) ---
$ws3Data = @(
    @("04:45:48", "05:43", "215A_LA PLATA", 58, "L6173"),
    @("05:24:16", "05:44", "215A_LA PLATA", 20, "L6173"),
    @("04:45:48", "06:08", "215A_LA PLATA", 83, "L6173"),
    @("05:55:02", "06:09", "215A_LA PLATA", 14, "L6173"),
    @("04:45:48", "06:32", "215C_LA PLATA", 107, "L6203"),
    @("06:25:28", "06:33", "215C_LA PLATA", 8, "L6203"),
    @("06:54:06", "07:00", "215B_LP-P MOR-1 Y 57", 6, "L6173"),
    @("07:17:59", "07:34", "215A_LA PLATA", 17, "L6173"),
    @("06:54:06", "07:35", "215A_LA PLATA", 41, "L6173"),
    @("06:25:28", "08:07", "215C_LA PLATA", 102, "L6203"),
    @("06:54:06", "08:13", "215C_LA PLATA", 79, "L6203"),
    @("07:17:59", "08:19", "215C_LA PLATA", 62, "L6203"),
    @("08:01:10", "08:22", "215C_LA PLATA", 21, "L6203"),
    @("07:48:05", "08:25", "215C_LA PLATA", 37, "L6203"),
    @("08:31:01", "08:33", "215C_LA PLATA", 2, "L6203"),
    @("07:17:59", "08:34", "215A_LA PLATA", 77, "L6173"),
    @("08:01:10", "08:35", "215A_LA PLATA", 34, "L6173"),
    @("08:31:01", "08:36", "215A_LA PLATA", 5, "L6173"),
    @("07:17:59", "09:08", "215D_LA PLATA", 111, "L6203"),
    @("08:47:26", "09:09", "215D_LA PLATA", 22, "L6203"),
    @("08:47:26", "10:03", "215B_LP-P MOR-40 Y 115", 76, "L6173")
)
for ($i = 0; $i -lt $ws3Data.Length; $i++) {
    $r = 6 + $i
    for ($c = 0; $c -lt 5; $c++) {
        $ws3.Cells.Item($r, $c + 1).Value = $ws3Data[$i][$c]
    }
}

Write-Host "Edit complete"